$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Record Number"
$ws.Range("B1").Value = "Latitude"
$ws.Range("C1").Value = "Longitude"
$ws.Range("D1").Value = "Roadway ID"
$ws.Range("E1").Value = "Accident Type"
$ws.Range("F1").Value = "Summary"

# --- Data rows ---
$data = @(
    @(1, 39.774537000000002, -86.175438999999997, "W. Michigan Street", "Other",   "V1 struck V2 while traveling west on W. Michigan Street"),
    @(2, 39.777228999999998, -86.178774000000004, "W. Walnut Street",   "Other",   "V1 struck pedestrian while they were crossing the street"),
    @(3, 39.778146000000000, -86.178729000000004, "Wishard Blvd.",      "Other",   "V1 rearended V2 near Walnut St. and Barnhill Dr."),
    @(4, 39.775764000000002, -86.186448999999996, "W. Michigan Street", "Weather", "V1 ran off the road due to ice on Michigan Street"),
    @(5, 39.778677999999999, -86.174892999999997, "University Blvd.",  "Other",   "V1 struck V2 while turning left from University Blvd. to Wishard Blvd."),
    @(6, 39.771799000000001, -86.178151999999997, "W. New York Street", "Other",   "V1 struck two pedestrians as they crossed the street from Barnhill Garage to Carroll Stadium"),
    @(7, 39.775593000000001, -86.183520000000001, "Eskenazi Ave.",      "Weather", "V1 slid off the road due severse weather conditions (ice)"),
    @(8, 39.777782000000002, -86.182849000000004, "Eskenazi Ave.",      "Other",   "V1 rearended a bus stopped at the bus stop while picking up passengers. One passenger fell from the bus steps during the incident and was taken to Eskenazi Hospital as a precaution")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 11.5

# --- Selection ---
$ws.Range("G12").Select()
